# "Generate Report for Handback"
# Fills in the Latest Target File / Latest Handback File / Latest Handback
# DateTime columns (I/J/K) for both language sheets, now that the handback
# round-trip has completed, and updates the Status column to reflect that
# the content is back in sync with en-US.

$wb = $excel.ActiveWorkbook

$statusOld = "In Translation"
$statusNew = "Handed back: in sync with en-US"

$url4dcaca83 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2350fc9bd3e5dd1ef68b79a38178f14d4284e907/e2e/4dcaca83-d5f8-4720-872f-ea3b50757ce3.md"
$url9a16f78c = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2350fc9bd3e5dd1ef68b79a38178f14d4284e907/e2e/9a16f78c-35fb-4931-9724-1e1c59fb8c4d.md"

$linkFontColor = 15570276  # BGR for RGB FF6495ED, matches the workbook's existing HyperLink style

# ---------------------------------------------------------------------
# Overview sheet: the Status text shown there is driven by the same
# shared "In Translation" text, so update it in place too.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusNew
$zhcn.Range("C3").Value = $statusNew

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $url4dcaca83, "", "", "4dcaca83-d5f8-4720-872f-ea3b50757ce3.md")
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = $linkFontColor
$zhcn.Range("J2").Value = "4dcaca83-d5f8-4720-872f-ea3b50757ce3.e01c17d3a8cd82c99bb046e69c175a8dd1cb5c75.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-23 10:24:41"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $url9a16f78c, "", "", "9a16f78c-35fb-4931-9724-1e1c59fb8c4d.md")
$zhcn.Range("I3").Font.Underline = $true
$zhcn.Range("I3").Font.Color = $linkFontColor
$zhcn.Range("J3").Value = "9a16f78c-35fb-4931-9724-1e1c59fb8c4d.f198b77d6c91c3daaaf6851d960e347918f82b5e.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-23 10:24:41"

$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusNew
$dede.Range("C3").Value = $statusNew

$dede.Hyperlinks.Add($dede.Range("I2"), $url4dcaca83, "", "", "4dcaca83-d5f8-4720-872f-ea3b50757ce3.md")
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = $linkFontColor
$dede.Range("J2").Value = "4dcaca83-d5f8-4720-872f-ea3b50757ce3.e01c17d3a8cd82c99bb046e69c175a8dd1cb5c75.de-de.xlf"
$dede.Range("K2").Value = "2016-08-23 10:24:48"

$dede.Hyperlinks.Add($dede.Range("I3"), $url9a16f78c, "", "", "9a16f78c-35fb-4931-9724-1e1c59fb8c4d.md")
$dede.Range("I3").Font.Underline = $true
$dede.Range("I3").Font.Color = $linkFontColor
$dede.Range("J3").Value = "9a16f78c-35fb-4931-9724-1e1c59fb8c4d.f198b77d6c91c3daaaf6851d960e347918f82b5e.de-de.xlf"
$dede.Range("K3").Value = "2016-08-23 10:24:48"

$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15
